# Finished all, need only logs
# Add the remaining Settings rows (email subject/body templates) and
# update the selected cell to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# New Name/Value pairs appended below the existing settings.
$ws.Cells.Item(10, 1).Value2 = "emailSubject"
$ws.Cells.Item(10, 2).Value2 = "Bitcoin Prices"

$ws.Cells.Item(11, 1).Value2 = "emailBodySuccess"
$ws.Cells.Item(11, 2).Value2 = "Dear user, here is the attachment from the run made at {0}."

$ws.Cells.Item(12, 1).Value2 = "emailBodyFailure"
$ws.Cells.Item(12, 2).Value2 = "Dear user, an error occurred."

# Update the active selection to B3, matching the saved view state.
$ws.Range("B3").Select() | Out-Null
